$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "type" column used to hold a single generic placeholder ("Type1") for
# every enemy. We now give each enemy its own concrete type name, starting
# with the existing three rows which all belong to the "Slime" enemy.
$ws.Range("B2:B4").Value = "Slime"

# Tune stats for the Slime rows (attackRange, attackSpeed, moveSpeed).
$ws.Range("E2").Value = 2
$ws.Range("E3").Value = 2
$ws.Range("E4").Value = 2

$ws.Range("I2").Value = 0.4
$ws.Range("I3").Value = 0.5
$ws.Range("I4").Value = 0.6

$ws.Range("J2").Value = 2

# Add the new TurtleShell (Normal rank) enemy on row 5.
$ws.Range("A5").Value = 10201
$ws.Range("B5").Value = "TurtleShell"
$ws.Range("C5").Value = "Normal"
$ws.Range("D5").Value = 30
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 3
$ws.Range("G5").Value = 0
$ws.Range("H5").Value = 3
$ws.Range("I5").Value = 0.4
$ws.Range("J5").Value = 2

# Let the type/rank columns resize to fit the new, longer text.
$ws.Columns.Item(2).AutoFit() | Out-Null
$ws.Columns.Item(3).AutoFit() | Out-Null

$ws.Range("I6").Select() | Out-Null
